$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text content updates ---
# (order matters for how new shared-string entries get appended)

# A30: repo renamed from fuck-university-physics-experiments to fuck-nku-physics-experiments
$ws.Range("A30").Value = "Posted on https://github.com/Axolyz/fuck-nku-physics-experiments."

# A3: clarify that pre-filled sample data should be replaced too
$ws.Range("A3").Value = "红色格子：填入你的实验数据，如本身自带数据请更改"

# A7: emphasize half-width colon requirement
$ws.Range("A7").Value = "请将度分秒数字以冒号(半角英文冒号!!)隔开输入，如输入：“54:30:00”，回车，单元格自动显示为 54°30′00″,输入角度时请写全度分秒"

# A8: add extra remark about this being the trickiest bit of the sheet
$ws.Range("A8").Value = "为表示角度，本表格利用了自带的日期格式，看到值突然变为日期是输入正确的体现, 这可能是整个灌水仓库里技术力最高的一个地方了"

# A32: drop the touhou quote, leaving the cell (and its protection style) empty
$ws.Range("A32").ClearContents() | Out-Null

# --- Formatting updates ---
# A7 and A8 get a new bold, teal-colored emphasis style. Start from A32's
# existing font (unlocked, 等线, family 3) so the new font clones that base
# instead of the plain default font, then tint it teal and bold it.
$teal = 12305721 # BGR encoding of RGB(0x39,0xC5,0xBB)

$ws.Range("A32").Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("A7").Font.Color = $teal
$ws.Range("A7").Font.Bold = $true

# A8 clones A7's now-finished format so no extra intermediate style gets
# materialized a second time.
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Sheet view updates ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A8").Select() | Out-Null
